$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.759.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.479.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.07%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.82"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.603"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.470.27"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.02%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.44"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.577"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "45.84"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.03%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.041.68"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.09%  "

$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "622.89"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -8.08%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.43"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -6.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.744.07"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.479.50"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.120"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.14"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.96"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.873"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.71"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -10.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.15"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.76"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.82%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.60"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -8.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -12.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.26"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -8.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.13"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -9.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.40"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -8.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.30"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -9.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.90"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "621.16"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.64"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.99%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.42%  "

$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -17.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.29"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0439"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.35%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.317.69"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.67%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.323"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.44"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0679"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -11.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.53"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -9.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.75"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.128"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.41"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.60"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +13.65%  "

